# Apply the Alvearie FHIR IG metadata refresh (v5.0.0 -> v6.0.0) to the
# "Metadata" sheet, and remove the duplicated "Contact" rows in favor of a
# single "Jurisdiction" row, shifting the rows below up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple value updates -------------------------------------------------
$ws.Range("B3").Value = "6.0.0"                               # Version
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"           # Date
$ws.Range("B9").Value = "Alvearie Team"                       # Publisher

# --- replace the two duplicated "Contact" rows (10 & 11) with a single
#     "Jurisdiction" row, then shift the remaining rows (Description,
#     Purpose, Copyright, Immutable) up by one and drop the now-empty last
#     row (15). -------------------------------------------------------------
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "The value set that indicates the place of service, based on standard values from Centers for Medicare and Medicaid Services (CMS)."

$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = $null

$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = $null

$ws.Range("A14").Value = "Immutable"
$ws.Range("B14").Value = "BooleanType[null]"

# Row 15 no longer exists in the refreshed sheet - delete it entirely so the
# used range shrinks back to A1:B14.
$ws.Range("A15:B15").Delete()
